# Edit workbook per commit: "added item: lesser healing potion ..."
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("items")

# E14: item gold value bump (rustysword-ish row) 15 -> 20
$ws.Range("E14").Value = 20

# New row 23: lesser healing potion (food item)
$ws.Range("A23").Value = "it_fo_lesserhealingpotion"
$ws.Range("B23").Value = "itd_lesserhealingpotion"
$ws.Range("C23").Value = 8
$ws.Range("D23").Value = "100, 100"
$ws.Range("E23").Value = 15
$ws.Range("G23").Value = 5
$ws.Range("T23").Value = 5
$ws.Range("U23").Value = "-15, -2"
$ws.Range("V23").Value = "20, 46"
$ws.Range("W23").Value = "100, 100, 50, 50"
$ws.Range("X23").Value = 1000

# Update the active selection to match the post-edit workbook state
$ws.Range("X25").Select() | Out-Null
